$wb = $excel.ActiveWorkbook

# --- Section_A sheet ---
$wsA = $wb.Worksheets.Item("Section_A")
$wsA.Range("B2").Value = "ELECTIVE_B6 [C404]"
$wsA.Range("C2").Value = "ELECTIVE_B7 [C304]"
$wsA.Range("D5").Value = "ELECTIVE_B6 [C404]"
$wsA.Range("E5").Value = "ELECTIVE_B7 [C304]"
$wsA.Range("C6").Value = "ELECTIVE_B6 (Tutorial) [C104]"
$wsA.Range("D6").Value = "ELECTIVE_B7 (Tutorial) [C304]"

# --- Section_B sheet ---
$wsB = $wb.Worksheets.Item("Section_B")
$wsB.Range("B2").Value = "ELECTIVE_B6 [C402]"
$wsB.Range("C2").Value = "ELECTIVE_B7 [C405]"
$wsB.Range("D5").Value = "ELECTIVE_B6 [C402]"
$wsB.Range("E5").Value = "ELECTIVE_B7 [C405]"
$wsB.Range("C6").Value = "ELECTIVE_B6 (Tutorial) [C402]"
$wsB.Range("D6").Value = "ELECTIVE_B7 (Tutorial) [C402]"

# --- Classroom_Utilization sheet ---
$wsU = $wb.Worksheets.Item("Classroom_Utilization")

# Row 6: C101
$wsU.Range("D6").Value = 0
$wsU.Range("E6").Value = 0
$wsU.Range("G6").Value = 0

# Row 9: C104
$wsU.Range("D9").Value = 1
$wsU.Range("E9").Value = 0.2
$wsU.Range("G9").Value = 2.5

# Row 15: C203
$wsU.Range("D15").Value = 0
$wsU.Range("E15").Value = 0
$wsU.Range("G15").Value = 0

# Row 23: C303
$wsU.Range("D23").Value = 0
$wsU.Range("E23").Value = 0
$wsU.Range("G23").Value = 0

# Row 24: C304
$wsU.Range("D24").Value = 4
$wsU.Range("E24").Value = 0.8
$wsU.Range("G24").Value = 10

# Row 29: C401
$wsU.Range("D29").Value = 0
$wsU.Range("E29").Value = 0
$wsU.Range("G29").Value = 0

# Row 30: C402
$wsU.Range("D30").Value = 5
$wsU.Range("E30").Value = 1
$wsU.Range("G30").Value = 12.5

# Row 33: C405
$wsU.Range("D33").Value = 3
$wsU.Range("E33").Value = 0.6
$wsU.Range("G33").Value = 7.5

# --- Classroom_Allocation sheet ---
$wsC = $wb.Worksheets.Item("Classroom_Allocation")

# The "Capacity" column (I) stores plain numeric-looking text (e.g. "96",
# "78") in the source workbook. Mark these cells as text first so Excel
# doesn't silently re-type them as numbers when the new value is assigned.
foreach ($cell in @("I2","I5","I8","I9","I11","I13")) {
    $wsC.Range($cell).NumberFormat = "@"
}

$wsC.Range("G2").Value = "C404"
$wsC.Range("I2").Value = "78"

$wsC.Range("G3").Value = "C304"

$wsC.Range("G4").Value = "C104"

$wsC.Range("G5").Value = "C404"
$wsC.Range("I5").Value = "78"

$wsC.Range("G6").Value = "C304"

$wsC.Range("G7").Value = "C304"

$wsC.Range("G8").Value = "C402"
$wsC.Range("I8").Value = "96"

$wsC.Range("G9").Value = "C405"
$wsC.Range("I9").Value = "78"

$wsC.Range("G10").Value = "C402"

$wsC.Range("G11").Value = "C402"
$wsC.Range("I11").Value = "96"

$wsC.Range("G12").Value = "C402"

$wsC.Range("G13").Value = "C405"
$wsC.Range("I13").Value = "78"
